$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 553.58826
$ws.Range("J41").Value = 737.6
$ws.Range("L41").Value = 737.6
$ws.Range("N41").Value = -1617.6

$ws.Range("H62").Value = 5607.154
$ws.Range("I62").Value = 4111.625
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 4111.625
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -3487.625
$ws.Range("N62").Value = -9248

$ws.Range("H65").Value = 5607.154
$ws.Range("I65").Value = 4111.625
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 20558.125
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -17438.125
$ws.Range("N65").Value = -46240

$ws.Range("H86").Value = 3980.5386
$ws.Range("J86").Value = 4604.9
$ws.Range("L86").Value = 4604.9
$ws.Range("N86").Value = -6850.9

$ws.Range("H89").Value = 3980.5386
$ws.Range("J89").Value = 4604.9
$ws.Range("L89").Value = 23024.5
$ws.Range("N89").Value = -34256.5

$ws.Range("H96").Value = 1912.5555
$ws.Range("I96").Value = 1151.625
$ws.Range("K96").Value = 3454.875
$ws.Range("M96").Value = -2081.875

$ws.Range("H113").Value = 3995
$ws.Range("I113").Value = 3995
$ws.Range("K113").Value = 3995
$ws.Range("M113").Value = -741

$ws.Range("H135").Value = 971.2857
$ws.Range("I135").Value = 731.4737
$ws.Range("K135").Value = 6583.263300000001
$ws.Range("M135").Value = -4048.263300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14755.196
$ws.Range("I32").Value = 2440.16
$ws.Range("J32").Value = 26596.576
$ws.Range("K32").Value = 2440.16
$ws.Range("L32").Value = 26596.576
$ws.Range("M32").Value = -2153.16
$ws.Range("N32").Value = -27170.576

$ws.Range("H63").Value = 5793.8237
$ws.Range("I63").Value = 4824.75
$ws.Range("K63").Value = 4824.75
$ws.Range("M63").Value = -4138.75

$ws.Range("H66").Value = 5793.8237
$ws.Range("I66").Value = 4824.75
$ws.Range("K66").Value = 24123.75
$ws.Range("M66").Value = -20691.75

$ws.Range("H122").Value = 436970.1
$ws.Range("I122").Value = 668097.6
$ws.Range("J122").Value = 3606
$ws.Range("K122").Value = 2004292.8
$ws.Range("L122").Value = 10818
$ws.Range("M122").Value = -2001842.8
$ws.Range("N122").Value = -15718

$ws.Range("H132").Value = 3076.5334
$ws.Range("I132").Value = 2558.9092
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 7676.7276
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -5146.7276
$ws.Range("N132").Value = -18560

$ws.Range("H140").Value = 105779.8
$ws.Range("J140").Value = 105779.8
$ws.Range("L140").Value = 105779.8
$ws.Range("N140").Value = -116139.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 979.9091
$ws.Range("I16").Value = 890
$ws.Range("J16").Value = 1384.5
$ws.Range("K16").Value = 890
$ws.Range("L16").Value = 1384.5
$ws.Range("M16").Value = -603
$ws.Range("N16").Value = -1958.5

$ws.Range("H86").Value = 13586.5
$ws.Range("I86").Value = 9449
$ws.Range("J86").Value = 14965.667
$ws.Range("K86").Value = 9449
$ws.Range("L86").Value = 14965.667
$ws.Range("M86").Value = -8326
$ws.Range("N86").Value = -17211.667

$ws.Range("H89").Value = 13586.5
$ws.Range("I89").Value = 9449
$ws.Range("J89").Value = 14965.667
$ws.Range("K89").Value = 47245
$ws.Range("L89").Value = 74828.33499999999
$ws.Range("M89").Value = -41629
$ws.Range("N89").Value = -86060.33499999999

$ws.Range("H93").Value = 2462.8333
$ws.Range("I93").Value = 2462.8333
$ws.Range("K93").Value = 2462.8333
$ws.Range("M93").Value = -590.8332999999998

$ws.Range("H105").Value = 2576.0625
$ws.Range("I105").Value = 691
$ws.Range("J105").Value = 4999.7144
$ws.Range("K105").Value = 691
$ws.Range("L105").Value = 4999.7144
$ws.Range("M105").Value = 1056
$ws.Range("N105").Value = -8493.7144

$ws.Range("H113").Value = 979.9091
$ws.Range("I113").Value = 890
$ws.Range("J113").Value = 1384.5
$ws.Range("K113").Value = 890
$ws.Range("L113").Value = 1384.5
$ws.Range("M113").Value = 1280
$ws.Range("N113").Value = -5724.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 195.33333
$ws.Range("I8").Value = 195.33333
$ws.Range("K8").Value = 585.99999
$ws.Range("M8").Value = -446.99999

$ws.Range("H113").Value = 951.6842
$ws.Range("J113").Value = 731
$ws.Range("L113").Value = 2193
$ws.Range("N113").Value = -6533

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31262

$ws.Range("H132").Value = 2605.1875
$ws.Range("I132").Value = 1733
$ws.Range("K132").Value = 5199
$ws.Range("M132").Value = -2669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2766.0833
$ws.Range("J46").Value = 3911.3333
$ws.Range("L46").Value = 3911.3333
$ws.Range("N46").Value = -4287.3333

$ws.Range("H55").Value = 330
$ws.Range("I55").Value = 330
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 330
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -157
$ws.Range("N55").ClearContents()

$ws.Range("H109").Value = 52000
$ws.Range("J109").Value = 52000
$ws.Range("L109").Value = 52000
$ws.Range("N109").Value = -54774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17116
$ws.Range("I45").Value = 25449.5
$ws.Range("J45").Value = 12949.25
$ws.Range("K45").Value = 25449.5
$ws.Range("L45").Value = 12949.25
$ws.Range("M45").Value = -24958.5
$ws.Range("N45").Value = -13931.25

$ws.Range("H81").Value = 1170.5
$ws.Range("I81").Value = 1170.5
$ws.Range("K81").Value = 2341
$ws.Range("M81").Value = -1280

$ws.Range("H84").Value = 1170.5
$ws.Range("I84").Value = 1170.5
$ws.Range("K84").Value = 11705
$ws.Range("M84").Value = -6401

$ws.Range("H107").Value = 825.0625
$ws.Range("I107").Value = 392.3846
$ws.Range("K107").Value = 1177.1538
$ws.Range("M107").Value = 742.8462

$ws.Range("H122").Value = 888.7143
$ws.Range("I122").Value = 888.7143
$ws.Range("K122").Value = 2666.1429
$ws.Range("M122").Value = -216.1428999999998

$ws.Range("H126").Value = 1882.5
$ws.Range("I126").Value = 730
$ws.Range("K126").Value = 2190
$ws.Range("M126").Value = 280

$ws.Range("H132").Value = 1370.5
$ws.Range("I132").Value = 1411.3
$ws.Range("J132").Value = 962.5
$ws.Range("K132").Value = 4233.9
$ws.Range("L132").Value = 2887.5
$ws.Range("M132").Value = -1703.9
$ws.Range("N132").Value = -7947.5
